# Update "countries & provincias Spain" dataset snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 16:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 503295
$ws.Range("C4").Value = 419
$ws.Range("D4").Value = 27335
$ws.Range("E4").Value = 457190
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 18770

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 122855
$ws.Range("C8").Value = 684
$ws.Range("E8").Value = 66206

# --- Row 19: Austria ---
$ws.Range("B19").Value = 13782
$ws.Range("C19").Value = 222
$ws.Range("E19").Value = 6841

# --- Row 23: Suecia ---
$ws.Range("F23").Value = 774

# --- Row 25: India ---
$ws.Range("B25").Value = 7997
$ws.Range("C25").Value = 397
$ws.Range("E25").Value = 6974

# --- Row 29: Australia ---
$ws.Range("B29").Value = 6303
$ws.Range("C29").Value = 65
$ws.Range("E29").Value = 2982

# --- Row 36: Pakistan ---
$ws.Range("B36").Value = 4970
$ws.Range("C36").Value = 275
$ws.Range("E36").Value = 4131

# --- Row 56: Argentina ---
$ws.Range("F56").Value = 115

# --- Rows 198-200: the underlying API re-ordered Nicaragua ahead of Cabo
# Verde / Santa Sede in the shared-string list, while the three rows (tied
# on total cases) kept their positions. Net effect: the row labels rotate
# (198 -> Nicaragua, 199 -> Cabo Verde, 200 -> Santa Sede) and Nicaragua's
# counters pick up one new/recovered case.
$ws.Range("A198").Value = "Nicaragua"
$ws.Range("B198").Value = 8
$ws.Range("C198").Value = 1
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Cabo Verde"
$ws.Range("B199").Value = 8
$ws.Range("C199").Value = 1
$ws.Range("D199").Value = 1
$ws.Range("E199").Value = 6
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

$ws.Range("A200").Value = "Santa Sede"
$ws.Range("B200").Value = 8
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 2
$ws.Range("E200").Value = 6
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0
